$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to text format so numeric-looking strings
# (e.g. "12546", "2000") are stored as shared strings, not numbers.
$numRng = $ws.Range("D2:E21")
$numRng.NumberFormat = "@"

$ws.Range("A2").Value = '에이치이엠파마(구.에이치이엠)'
$ws.Range("B2").Value = '2024.08.26~08.30'
$ws.Range("C2").Value = '18,000~21,000'
$ws.Range("D2").Value = '-'
$ws.Range("E2").Value = '12546'
$ws.Range("F2").Value = '신한투자증권'

$ws.Range("A3").Value = '미래에셋비전스팩7호'
$ws.Range("B3").Value = '2024.08.26~08.27'
$ws.Range("C3").Value = '2,000~2,000'
$ws.Range("D3").Value = '-'
$ws.Range("E3").Value = '-'
$ws.Range("F3").Value = '미래에셋증권'

$ws.Range("A4").Value = '아이언디바이스'
$ws.Range("B4").Value = '2024.08.19~08.23'
$ws.Range("C4").Value = '4,900~5,700'
$ws.Range("D4").Value = '-'
$ws.Range("E4").Value = '14700'
$ws.Range("F4").Value = '대신증권'

$ws.Range("A5").Value = '키움스팩9호'
$ws.Range("B5").Value = '2024.08.13~08.14'
$ws.Range("C5").Value = '2,000~2,000'
$ws.Range("D5").Value = '-'
$ws.Range("E5").Value = '6000'
$ws.Range("F5").Value = '키움증권'

$ws.Range("A6").Value = '아이스크림미디어(구.시공미디어)'
$ws.Range("B6").Value = '2024.08.09~08.16'
$ws.Range("C6").Value = '32,000~40,200'
$ws.Range("D6").Value = '-'
$ws.Range("E6").Value = '78720'
$ws.Range("F6").Value = '삼성증권'

$ws.Range("A7").Value = '이엔셀'
$ws.Range("B7").Value = '2024.08.02~08.08'
$ws.Range("C7").Value = '13,600~15,300'
$ws.Range("D7").Value = '-'
$ws.Range("E7").Value = '21308'
$ws.Range("F7").Value = 'NH투자증권'

$ws.Range("A8").Value = '엠83'
$ws.Range("B8").Value = '2024.08.01~08.07'
$ws.Range("C8").Value = '11,000~13,000'
$ws.Range("D8").Value = '-'
$ws.Range("E8").Value = '16500'
$ws.Range("F8").Value = '신영증권,유진투자증권'

$ws.Range("A9").Value = '티디에스팜'
$ws.Range("B9").Value = '2024.07.31~08.06'
$ws.Range("C9").Value = '9,500~10,700'
$ws.Range("D9").Value = '-'
$ws.Range("E9").Value = '9500'
$ws.Range("F9").Value = '한국투자증권'

$ws.Range("A10").Value = '케이쓰리아이'
$ws.Range("B10").Value = '2024.07.30~08.05'
$ws.Range("C10").Value = '12,500~15,500'
$ws.Range("D10").Value = '-'
$ws.Range("E10").Value = '17500'
$ws.Range("F10").Value = '하나증권'

$ws.Range("A11").Value = '전진건설로봇(구.전진중공업)(유가)'
$ws.Range("B11").Value = '2024.07.30~08.05'
$ws.Range("C11").Value = '13,800~15,700'
$ws.Range("D11").Value = '-'
$ws.Range("E11").Value = '42471'
$ws.Range("F11").Value = '미래에셋증권'

$ws.Range("A12").Value = '교보스팩16호'
$ws.Range("B12").Value = '2024.07.29~07.30'
$ws.Range("C12").Value = '2,000~2,000'
$ws.Range("D12").Value = '2000'
$ws.Range("E12").Value = '11600'
$ws.Range("F12").Value = '교보증권'

$ws.Range("A13").Value = '넥스트바이오메디컬'
$ws.Range("B13").Value = '2024.07.29~08.02'
$ws.Range("C13").Value = '24,000~29,000'
$ws.Range("D13").Value = '-'
$ws.Range("E13").Value = '24000'
$ws.Range("F13").Value = '한국투자증권'

$ws.Range("A14").Value = '유라클'
$ws.Range("B14").Value = '2024.07.29~08.02'
$ws.Range("C14").Value = '18,000~21,000'
$ws.Range("D14").Value = '-'
$ws.Range("E14").Value = '13518'
$ws.Range("F14").Value = '키움증권'

$ws.Range("A15").Value = '뱅크웨어글로벌'
$ws.Range("B15").Value = '2024.07.23~07.29'
$ws.Range("C15").Value = '16,000~19,000'
$ws.Range("D15").Value = '16000'
$ws.Range("E15").Value = '22400'
$ws.Range("F15").Value = '미래에셋증권'

$ws.Range("A16").Value = '아이빔테크놀로지'
$ws.Range("B16").Value = '2024.07.15~07.19'
$ws.Range("C16").Value = '7,300~8,500'
$ws.Range("D16").Value = '10000'
$ws.Range("E16").Value = '16308'
$ws.Range("F16").Value = '삼성증권'

$ws.Range("A17").Value = '피앤에스미캐닉스'
$ws.Range("B17").Value = '2024.07.11~07.17'
$ws.Range("C17").Value = '14,000~17,000'
$ws.Range("D17").Value = '22000'
$ws.Range("E17").Value = '18900'
$ws.Range("F17").Value = '키움증권'

$ws.Range("A18").Value = 'NH스팩31호'
$ws.Range("B18").Value = '2024.07.09~07.10'
$ws.Range("C18").Value = '2,000~2,000'
$ws.Range("D18").Value = '2000'
$ws.Range("E18").Value = '12000'
$ws.Range("F18").Value = 'NH투자증권'

$ws.Range("A19").Value = 'SK증권스팩13호'
$ws.Range("B19").Value = '2024.07.09~07.10'
$ws.Range("C19").Value = '2,000~2,000'
$ws.Range("D19").Value = '2000'
$ws.Range("E19").Value = '8000'
$ws.Range("F19").Value = 'SK증권'

$ws.Range("A20").Value = '산일전기(유가)'
$ws.Range("B20").Value = '2024.07.09~07.15'
$ws.Range("C20").Value = '24,000~30,000'
$ws.Range("D20").Value = '35000'
$ws.Range("E20").Value = '182400'
$ws.Range("F20").Value = '미래에셋증권,삼성증권'

$ws.Range("A21").Value = '이베스트스팩6호'
$ws.Range("B21").Value = '2024.06.27~06.28'
$ws.Range("C21").Value = '2,000~2,000'
$ws.Range("D21").Value = '2000'
$ws.Range("E21").Value = '8000'
$ws.Range("F21").Value = '엘에스증권'

# Restore default (unformatted) style on D2:E21 while keeping the
# values/types (ClearFormats resets style but preserves cell content).
$numRng.ClearFormats()
